# Game Testing Log - add win screen test cases and a few more rows,
# and append a trailing period to each "Expected result" that already
# existed, writing the updated text into the previously-empty "Actual
# Result" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the missing trailing-period copies in column D for the three
#     existing rows (moving left / moving Right / Jumping Up) ---
$ws.Range("D2").Value = "playerXpos decreases and player image moves left on canvas."
$ws.Range("D4").Value = "PlayerXpos Increases and PlayerImage moves right on Canvas."
$ws.Range("D6").Value = "PlayerYpos Increases  and player image also shifts up but then needs to wait 1 sec to go up the canvas do again."

# --- Row 8: Contact with box ---
$ws.Range("A8").Value = "Contact with box"
$ws.Range("B8").Value = "any arrow key"
$ws.Range("C8").Value = "Box should cause the players pos to move a different way to create the illusion of a collision"
$ws.Range("D8").Value = "Box should cause the players pos to move a different way to create the illusion of a collision."

# --- Row 10: Contact with Heart (the new win screen) ---
$ws.Range("A10").Value = "Contact with Heart"
$ws.Range("B10").Value = "any arrow key"
$ws.Range("C10").Value = 'When in contact with Heart area the Winning screen coshould me up displaying the name and "Congrats"'
$ws.Range("D10").Value = 'When in contact with Heart area the Winning screen coshould me up displaying the name and "Congrats".'

# --- Row 12: Gravity ---
$ws.Range("A12").Value = "Gravity"
$ws.Range("B12").Value = "Not moving"
$ws.Range("C12").Value = 'When PLAYER_SIZE past the canvas width or height black death screen should come up displaying player name and "Nice try"'
$ws.Range("D12").Value = 'When PLAYER_SIZE past the canvas width or height black death screen should come up displaying player name and "Nice try".'

# --- Row 14: Spawn ---
$ws.Range("A14").Value = "Spawn"
$ws.Range("B14").Value = "Start game"
$ws.Range("C14").Value = "Player starting pos should be to the right of the map at the start"
$ws.Range("D14").Value = "Player starting pos should be to the right of the map at the start."

# --- Match row heights of the newly-populated rows to the other
#     populated rows in the sheet (wrapped multi-line text) ---
$ws.Range("A8:D8").RowHeight = 33
$ws.Range("A10:D10").RowHeight = 49.5
$ws.Range("A12:D12").RowHeight = 49.5
$ws.Range("A14:D14").RowHeight = 33

# --- Update active selection to match the source workbook ---
$ws.Range("E6").Select() | Out-Null
